$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column N (year 2022) values, mirroring the existing style used by the
# corresponding cell in column M (or, where that already carries the "0.0"
# number format, reusing that same style) plus the "0.0" number format.

$ws.Range("N2").Value = ""
$ws.Range("N2").NumberFormat = "General"

$ws.Range("N3").Value = 2022

$ws.Range("N4").Value = 9.224468514531754
$ws.Range("N4").NumberFormat = "0.0"

$ws.Range("N5").Value = 4.6068543125097872
$ws.Range("N5").NumberFormat = "0.0"

$ws.Range("N6").Value = 13.543910285971602
$ws.Range("N6").NumberFormat = "0.0"

$ws.Range("N7").Value = 24.703327617190443
$ws.Range("N7").NumberFormat = "0.0"

$ws.Range("N8").Value = 28.608474183838851
$ws.Range("N8").NumberFormat = "0.0"

$ws.Range("N9").Value = 20.904451081350146
$ws.Range("N9").NumberFormat = "0.0"

$ws.Range("N10").Value = 26.720095429750884
$ws.Range("N10").NumberFormat = "0.0"

$ws.Range("N11").Value = 27.704327204727914
$ws.Range("N11").NumberFormat = "0.0"

$ws.Range("N12").Value = 25.731792255708452
$ws.Range("N12").NumberFormat = "0.0"

# Match fonts/borders per row to mirror column M's look for the new column.
$ws.Range("N2").Font.Name = "Times New Roman"
$ws.Range("N2").Font.Size = 9
$ws.Range("N2").Font.Bold = $false

$ws.Range("N3").Font.Name = "Times New Roman"
$ws.Range("N3").Font.Size = 9
$ws.Range("N3").Font.Bold = $true

foreach ($r in 4,7,10) {
    $ws.Range("N$r").Font.Name = "Times New Roman"
    $ws.Range("N$r").Font.Size = 9
    $ws.Range("N$r").Font.Bold = $true
}

foreach ($r in 5,6,8,9,11,12) {
    $ws.Range("N$r").Font.Name = "Times New Roman"
    $ws.Range("N$r").Font.Size = 9
    $ws.Range("N$r").Font.Bold = $false
}

foreach ($r in 2,3,4,5,6,7,8,9,10,11,12) {
    $ws.Range("N$r").VerticalAlignment = -4108
}

$ws.Range("N12").Borders.Item(9).LineStyle = -4138
$ws.Range("N12").Borders.Item(9).Weight = -4138

# Move the active selection the way the recorded session left it.
$ws.Range("Q5").Select()
